$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 2646
$ws.Range("J3").Value = 2716
$ws.Range("F4").Value = 1884
$ws.Range("J4").Value = 617
$ws.Range("J5").Value = 213
$ws.Range("J6").Value = 3349
$ws.Range("F7").Value = 24074
$ws.Range("J7").Value = 9541

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J2").Value = 74
$ws.Range("J6").Value = 92
$ws.Range("F7").Value = 827
$ws.Range("J7").Value = 287
$ws.Range("J8").Value = 597
$ws.Range("J14").Value = 33
$ws.Range("J15").Value = 114
$ws.Range("J19").Value = 306
$ws.Range("J20").Value = 196
$ws.Range("J24").Value = 29
$ws.Range("J29").Value = 556
$ws.Range("J30").Value = 37
$ws.Range("J31").Value = 72
$ws.Range("J33").Value = 390
$ws.Range("J36").Value = 143
$ws.Range("J37").Value = 324
$ws.Range("J41").Value = 64
$ws.Range("J42").Value = 372
$ws.Range("J44").Value = 78
$ws.Range("J46").Value = 32
$ws.Range("J48").Value = 95
$ws.Range("J49").Value = 60
$ws.Range("J52").Value = 240
$ws.Range("J54").Value = 187
$ws.Range("J55").Value = 113
$ws.Range("J57").Value = 44
$ws.Range("J63").Value = 48
$ws.Range("J64").Value = 63
$ws.Range("J65").Value = 247
$ws.Range("J67").Value = 344
$ws.Range("J73").Value = 86
$ws.Range("J76").Value = 133
$ws.Range("J77").Value = 80
$ws.Range("J79").Value = 287
$ws.Range("J85").Value = 440
$ws.Range("J89").Value = 99
$ws.Range("J91").Value = 109
$ws.Range("J95").Value = 149
$ws.Range("J96").Value = 113
$ws.Range("J97").Value = 60
$ws.Range("J98").Value = 57
$ws.Range("J99").Value = 132
$ws.Range("F101").Value = 24074
$ws.Range("J101").Value = 9541

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 74

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J3").Value = 26
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 96
$ws.Range("J3").Value = 82
$ws.Range("F4").Value = 55
$ws.Range("J6").Value = 100
$ws.Range("F7").Value = 827
$ws.Range("J7").Value = 287

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J2").Value = 188
$ws.Range("J4").Value = 29
$ws.Range("J6").Value = 170
$ws.Range("J7").Value = 597

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 33
$ws.Range("J7").Value = 114

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 85
$ws.Range("J7").Value = 306

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 63
$ws.Range("J3").Value = 59
$ws.Range("J7").Value = 196

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 162
$ws.Range("J3").Value = 185
$ws.Range("J6").Value = 154
$ws.Range("J7").Value = 556

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("J6").Value = 10
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("J2").Value = 29
$ws.Range("J3").Value = 17
$ws.Range("J7").Value = 72

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J3").Value = 119
$ws.Range("J6").Value = 132
$ws.Range("J7").Value = 390

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J2").Value = 52
$ws.Range("J3").Value = 37
$ws.Range("J7").Value = 143

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 114
$ws.Range("J5").Value = 14
$ws.Range("J6").Value = 91
$ws.Range("J7").Value = 324

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 64

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 75
$ws.Range("J3").Value = 83
$ws.Range("J6").Value = 188
$ws.Range("J7").Value = 372

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J3").Value = 21
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 78

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("J3").Value = 14
$ws.Range("J6").Value = 46
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 60

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J2").Value = 55
$ws.Range("J4").Value = 13
$ws.Range("J7").Value = 240

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J3").Value = 38
$ws.Range("J7").Value = 187

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J3").Value = 21
$ws.Range("J7").Value = 113

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("J6").Value = 18
$ws.Range("J7").Value = 44

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J6").Value = 89
$ws.Range("J7").Value = 247

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 137
$ws.Range("J6").Value = 103
$ws.Range("J7").Value = 344

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J3").Value = 26
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J6").Value = 74
$ws.Range("J7").Value = 133

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("J2").Value = 25
$ws.Range("J3").Value = 27
$ws.Range("J5").Value = 4
$ws.Range("J7").Value = 80

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 81
$ws.Range("J3").Value = 106
$ws.Range("J6").Value = 79
$ws.Range("J7").Value = 287

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 105
$ws.Range("J6").Value = 127
$ws.Range("J7").Value = 440

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J2").Value = 33
$ws.Range("J7").Value = 99

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 31
$ws.Range("J3").Value = 48
$ws.Range("J7").Value = 109

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 43
$ws.Range("J6").Value = 41
$ws.Range("J7").Value = 149

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J2").Value = 37
$ws.Range("J3").Value = 33
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 113

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 60

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J6").Value = 34
$ws.Range("J7").Value = 132
